$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells keep their original text representation
# (values like "596.21" or "0.0000183" would otherwise be auto-converted
# to numbers by Excel's type inference): force text format before writing
# the value, then restore the default "Normal" style so the cell's
# appearance/format stays identical to the rest of the sheet.
$priceCells = @("D2","D3","D5","D6","D9","D10","D13","D14","D15","D16","D17","D18","D19","D20","D21","D23","D26","D27","D28","D29","D30","D31","D32","D36","D37","D38","D46","D47","D48","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "71.799.23"
$ws.Range("E2").Value = "  +0.95%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.688.13"
$ws.Range("E3").Value = "  +2.67%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "596.21"
$ws.Range("E5").Value = "  -1.46%  "

# Row 6 - Solana
$ws.Range("D6").Value = "174.85"
$ws.Range("E6").Value = "  -2.20%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.67%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.686.66"
$ws.Range("E9").Value = "  +2.52%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +1.31%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +2.31%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.22%  "

# Row 13 - Toncoin
$ws.Range("D13").Value = "4.99"

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.177.20"
$ws.Range("E14").Value = "  +1.20%  "

# Row 15 - was ShibaInu, now WrappedBTC
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "71.902.85"
$ws.Range("E15").Value = "  +1.13%  "

# Row 16 - was WrappedBTC, now ShibaInu
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000183"
$ws.Range("E16").Value = "  -1.46%  "

# Row 17 - Avalanche
$ws.Range("D17").Value = "26.11"
$ws.Range("E17").Value = "  -2.09%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.643.70"
$ws.Range("E18").Value = "  +0.77%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "12.05"
$ws.Range("E19").Value = "  +5.05%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "8.00"
$ws.Range("E20").Value = "  +1.45%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "369.69"
$ws.Range("E21").Value = "  -2.95%  "

# Row 22 - Polkadot
$ws.Range("E22").Value = "  -0.38%  "

# Row 23 - SuiNetwork
$ws.Range("D23").Value = "2.01"
$ws.Range("E23").Value = "  +1.01%  "

# Row 24 - Litecoin
$ws.Range("E24").Value = "  -0.91%  "

# Row 26 - NEARProtocol
$ws.Range("D26").Value = "4.29"
$ws.Range("E26").Value = "  -3.59%  "

# Row 27 - Aptos
$ws.Range("D27").Value = "9.80"
$ws.Range("E27").Value = "  -2.54%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "2.721.76"
$ws.Range("E28").Value = "  -1.22%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.26%  "

# Row 30 - PEPE
$ws.Range("D30").Value = "0.0₃0946"
$ws.Range("E30").Value = "  -0.95%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "8.04"
$ws.Range("E31").Value = "  -0.16%  "

# Row 32 - Bittensor
$ws.Range("D32").Value = "503.77"
$ws.Range("E32").Value = "  -7.86%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  -3.95%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  -1.05%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  -0.20%  "

# Row 36 - Monero
$ws.Range("D36").Value = "162.75"
$ws.Range("E36").Value = "  -2.12%  "

# Row 37 - EthereumClassic
$ws.Range("D37").Value = "19.44"
$ws.Range("E37").Value = "  +1.19%  "

# Row 38 - WhiteBITCoin
$ws.Range("D38").Value = "19.07"
$ws.Range("E38").Value = "  -0.45%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -2.69%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  -6.15%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -4.37%  "

# Row 42 - USDe
$ws.Range("E42").Value = "  -0.05%  "

# Row 43 - RenderToken
$ws.Range("E43").Value = "  -1.15%  "

# Row 44 - dogwifhat
$ws.Range("E44").Value = "  -2.55%  "

# Row 45 - PolygonEcosystemToken
$ws.Range("E45").Value = "  -0.42%  "

# Row 46 - was OKB, now Aave
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "156.14"
$ws.Range("E46").Value = "  +2.12%  "

# Row 47 - was Aave, now OKB
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "39.18"
$ws.Range("E47").Value = "  -2.04%  "

# Row 48 - Filecoin
$ws.Range("D48").Value = "3.70"
$ws.Range("E48").Value = "  +1.78%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  +2.42%  "

# Row 50 - Optimism
$ws.Range("E50").Value = "  +2.74%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "0.0763"
$ws.Range("E51").Value = "  +0.70%  "

# Restore the default "Normal" cell style on the price cells so only the
# text content changed (the temporary "@" text format was only needed to
# stop Excel from re-interpreting the strings as numbers on entry).
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
